$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.283252000808716
$ws.Range("B1").Value = 1.425832629203796
$ws.Range("C1").Value = 6.82899808883667
$ws.Range("D1").Value = 1.993505716323853
$ws.Range("E1").Value = 0.8935665488243103
